$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The simulation re-run collapsed every winner that used to be one of
# Zemmour / Hidalgo / Pecresse / Jadot (old shared-string slots 12-15)
# into either Montebourg or Philipot across rows 4..12 (columns B..AQ).
# B4 additionally lost its value entirely (no winner recorded there).
#
# Row patterns below cover columns B..AQ (worksheet column index 2..43):
#   "M" -> Montebourg
#   "P" -> Philipot
#   "." -> clear the cell (no value)
$rowPatterns = @{
    4  = ".PMMMMMPMMMMMPMMMMMPMMMMMPMMMMMPMMMMMPMMMM"
    5  = "MPMMMMMPMMMMMPMMMMMMMMMMMPMMMMMPMMMMMPMMMM"
    6  = "MPMMMMMPMMMMMPMMMMMPMMMMMPMMMMMPMMMMMPMMMM"
    7  = "MPMMMMMMMMMMMMMMMMMPMMMMMMMMMMMPMMMMPPMMMM"
    8  = "MPMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMM"
    9  = "MPMMMMMPMMMMMPMMMMMPMMMMMPMMMMMPMMMMMPMMMM"
    10 = "MMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMM"
    11 = "MMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMM"
    12 = "MPMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMMM"
}

foreach ($row in $rowPatterns.Keys) {
    $pattern = $rowPatterns[$row]
    for ($i = 0; $i -lt $pattern.Length; $i++) {
        $col = $i + 2   # column B is worksheet column index 2
        $code = $pattern.Substring($i, 1)
        $cell = $ws.Cells.Item($row, $col)
        if ($code -eq ".") {
            $cell.ClearContents()
        } elseif ($code -eq "M") {
            $cell.Value = "Montebourg"
        } elseif ($code -eq "P") {
            $cell.Value = "Philipot"
        }
    }
}
